$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Remove the old subtitle row (row 2: "pompes)", "Hiver", "Eté", "Année", unit labels)
#    This shifts the data rows (old 3-7) up to become rows 2-6.
$ws.Rows("2:2").Delete()

# 2) Rewrite the header row (row 1) with the new column layout:
#    A: idx, B: idx2, C: Name, D: Date Start, E: Date End,
#    F: (m3/s), G: (MW1), H: (MW2), I: (GWh) Winter, J: (GWh) Summer, K: (GWh) Year
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 pick up the (new) font-only style used for header units.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# 3) Update the selection to match the target workbook
$ws.Range("A2:K2").Select()

Write-Output "done"
